{"js": "// 1. Touch the section's even/first-page headers & footers so Word mints\n//    header1.xml (even), header2.xml (default/primary), header3.xml (first)\n//    and footer1.xml (even), footer2.xml (default/primary), footer3.xml (first),\n//    wiring the six header/footer references into sectPr.\nconst sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\n\nconst section = sections.items[0];\nconst evenHeader = section.getHeader(\"EvenPages\");\nconst primaryHeader = section.getHeader(\"Primary\");\nconst firstHeader = section.getHeader(\"FirstPage\");\nconst evenFooter = section.getFooter(\"EvenPages\");\nconst primaryFooter = section.getFooter(\"Primary\");\nconst firstFooter = section.getFooter(\"FirstPage\");\nawait context.sync();\n\n// 2. Rewrite the primary header text (now persisted as header2.xml) so it\n//    reads \"CS164 \u2013 Worksheet Week 2 \u2013 Finding Errors\" as separate runs,\n//    and re-append the _GoBack bookmark at the end of that paragraph.\nconst headerOoxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\"><pkg:xmlData>\n<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\"><Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/></Relationships>\n</pkg:xmlData></pkg:part>\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>\n<w:p>\n<w:pPr><w:pStyle w:val=\"Header\"/></w:pPr>\n<w:r><w:t xml:space=\"preserve\">CS164 \u2013 Worksheet Week </w:t></w:r>\n<w:r><w:t>2</w:t></w:r>\n<w:r><w:t xml:space=\"preserve\"> \u2013 </w:t></w:r>\n<w:r><w:t xml:space=\"preserve\">Finding </w:t></w:r>\n<w:r><w:t>Errors</w:t></w:r>\n<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>\n<w:bookmarkEnd w:id=\"0\"/>\n</w:p>\n</w:body></w:document>\n</pkg:xmlData></pkg:part>\n</pkg:package>`;\nprimaryHeader.insertOoxml(headerOoxml, Word.InsertLocation.replace);\nawait context.sync();\n\n// 3. Remove the stray _GoBack bookmark that used to live in the body.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$sec = $d.Sections(1)\n\n# Turning on odd/even headers is what makes Word mint the even-page header\n# (header1.xml), keep the primary header as the \"default\" header\n# (header2.xml), and also wires up the first-page header (header3.xml) plus\n# the matching even/default/first footers (footer1-3.xml) into sectPr.\n$sec.PageSetup.OddAndEvenPagesHeaderFooter = 1\n\n$hPrimary = $sec.Headers(1)\n$hFirst   = $sec.Headers(2)\n$hEven    = $sec.Headers(3)\n$fPrimary = $sec.Footers(1)\n$fFirst   = $sec.Footers(2)\n$fEven    = $sec.Footers(3)\n\n# Touch the non-primary headers/footers so their parts get created empty.\n$hFirst.Range.Text = \"\"\n$hEven.Range.Text = \"\"\n$fPrimary.Range.Text = \"\"\n$fFirst.Range.Text = \"\"\n$fEven.Range.Text = \"\"\n\n# Rewrite the primary header (persisted as header2.xml) with the new title,\n# split across runs, and re-append the _GoBack bookmark at the end of the\n# paragraph.\n$ooxml = @'\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\"><pkg:xmlData>\n<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\"><Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/></Relationships>\n</pkg:xmlData></pkg:part>\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>\n<w:p>\n<w:pPr><w:pStyle w:val=\"Header\"/></w:pPr>\n<w:r><w:t xml:space=\"preserve\">CS164 \u2013 Worksheet Week </w:t></w:r>\n<w:r><w:t>2</w:t></w:r>\n<w:r><w:t xml:space=\"preserve\"> \u2013 </w:t></w:r>\n<w:r><w:t xml:space=\"preserve\">Finding </w:t></w:r>\n<w:r><w:t>Errors</w:t></w:r>\n<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>\n<w:bookmarkEnd w:id=\"0\"/>\n</w:p>\n</w:body></w:document>\n</pkg:xmlData></pkg:part>\n</pkg:package>\n'@\n$hPrimary.Range.InsertXML($ooxml)\n\n# Remove the stray _GoBack bookmark that used to live in the body.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n"}
